$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "111×2=222" "316×5=1580"
Replace-Text "732×4=2928" "343×7=2401"
Replace-Text "220×9=1980" "861×8=6888"
Replace-Text "760×7=5320" "623×5=3115"
Replace-Text "461×7=3227" "866×2=1732"
Replace-Text "462×8=3696" "903×9=8127"
Replace-Text "864×4=3456" "303×8=2424"
Replace-Text "614×8=4912" "615×9=5535"
Replace-Text "576×7=4032" "132×5=660"
Replace-Text "327×6=1962" "740×2=1480"
Replace-Text "499×2=998" "925×3=2775"
Replace-Text "213×6=1278" "844×8=6752"
Replace-Text "228×7=1596" "525×2=1050"
Replace-Text "777×3=2331" "801×6=4806"
Replace-Text "544×6=3264" "566×8=4528"
Replace-Text "763×6=4578" "538×4=2152"
Replace-Text "964×9=8676" "514×4=2056"
Replace-Text "717×2=1434" "724×6=4344"
Replace-Text "184×5=920" "212×2=424"
Replace-Text "702×7=4914" "109×9=981"
Replace-Text "587×6=3522" "733×9=6597"
Replace-Text "102×7=714" "548×5=2740"
Replace-Text "780×6=4680" "621×6=3726"
Replace-Text "989×5=4945" "428×4=1712"
Replace-Text "691×2=1382" "163×8=1304"
